$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
if ($lastRow -lt 2) { $lastRow = 381 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $val = $cell.Value2
    if ($null -eq $val) { continue }
    $text = [string]$val
    if ($text.Length -eq 0) { continue }

    if ($text.StartsWith("m.")) {
        $newText = "miasto " + $text.Substring(2)
    } else {
        $newText = "powiat " + $text
    }

    $cell.Value2 = $newText
}
